$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A-T), reflecting the NATMI
# re-run with 3 resampling replicates ("Natmi following Dr Hou advice").
# Column order: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity,
# T Edge total expression derived specificity.
$data = @(
    @('ECs', 'Cp', 'Slc40a1', 'ECs', 3, 1, 7.02243, 21.06729, 0.09918483354594591, 0.0991848335459459, 2, 0.6666666666666666, 0.1166653333333333, 0.349996, 0.003627699309917676, 0.003627699309917676, 0.8192741367599999, 7.373467230839999, 0.0003598127522089276, 0.0003598127522089276),
    @('ECs', 'Cp', 'Slc40a1', 'FAPs', 3, 1, 7.02243, 21.06729, 0.09918483354594591, 0.0991848335459459, 3, 1, 12.00209833333333, 36.00629499999999, 0.3732042981182421, 0.3732042981182421, 84.28389539894998, 758.5550585905498, 0.03701620618748942, 0.03701620618748942),
    @('ECs', 'Cp', 'Slc40a1', 'M2', 3, 1, 7.02243, 21.06729, 0.09918483354594591, 0.0991848335459459, 3, 1, 14.050981, 42.15294300000001, 0.4369141425390552, 0.4369141425390553, 98.67203050383002, 888.0482745344701, 0.04333525650160588, 0.04333525650160588),
    @('ECs', 'Cp', 'Slc40a1', 'sCs', 3, 1, 7.02243, 21.06729, 0.09918483354594591, 0.0991848335459459, 3, 1, 5.989848333333334, 17.969545, 0.186253860032785, 0.186253860032785, 42.06329063145, 378.56961568305, 0.01847355810464169, 0.01847355810464169),
    @('FAPs', 'Cp', 'Slc40a1', 'ECs', 3, 1, 37.70032666666666, 113.10098, 0.5324795868468776, 0.5324795868468776, 2, 0.6666666666666666, 0.1166653333333333, 0.349996, 0.003627699309917676, 0.003627699309917676, 4.398321177342222, 39.58489059607999, 0.001931675829749667, 0.001931675829749667),
    @('FAPs', 'Cp', 'Slc40a1', 'FAPs', 3, 1, 37.70032666666666, 113.10098, 0.5324795868468776, 0.5324795868468776, 3, 1, 12.00209833333333, 36.00629499999999, 0.3732042981182421, 0.3732042981182421, 452.4830278521221, 4072.347250669099, 0.1987236704714805, 0.1987236704714805),
    @('FAPs', 'Cp', 'Slc40a1', 'M2', 3, 1, 37.70032666666666, 113.10098, 0.5324795868468776, 0.5324795868468776, 3, 1, 14.050981, 42.15294300000001, 0.4369141425390552, 0.4369141425390553, 529.7265736871267, 4767.53916318414, 0.2326478621067539, 0.2326478621067539),
    @('FAPs', 'Cp', 'Slc40a1', 'sCs', 3, 1, 37.70032666666666, 113.10098, 0.5324795868468776, 0.5324795868468776, 3, 1, 5.989848333333334, 17.969545, 0.186253860032785, 0.186253860032785, 225.8192388504555, 2032.3731496541, 0.09917637843889354, 0.09917637843889354),
    @('M2', 'Cp', 'Slc40a1', 'ECs', 3, 1, 4.606340666666667, 13.819022, 0.06505997671450693, 0.06505997671450692, 2, 0.6666666666666666, 0.1166653333333333, 0.349996, 0.003627699309917676, 0.003627699309917676, 0.5374002693235556, 4.836602423912, 0.0002360180326304769, 0.0002360180326304769),
    @('M2', 'Cp', 'Slc40a1', 'FAPs', 3, 1, 4.606340666666667, 13.819022, 0.06505997671450693, 0.06505997671450692, 3, 1, 12.00209833333333, 36.00629499999999, 0.3732042981182421, 0.3732042981182421, 55.28575363816555, 497.57178274349, 0.02428066294532673, 0.02428066294532673),
    @('M2', 'Cp', 'Slc40a1', 'M2', 3, 1, 4.606340666666667, 13.819022, 0.06505997671450693, 0.06505997671450692, 3, 1, 14.050981, 42.15294300000001, 0.4369141425390552, 0.4369141425390553, 64.72360518686068, 582.5124466817462, 0.02842562393982969, 0.02842562393982969),
    @('M2', 'Cp', 'Slc40a1', 'sCs', 3, 1, 4.606340666666667, 13.819022, 0.06505997671450693, 0.06505997671450692, 3, 1, 5.989848333333334, 17.969545, 0.186253860032785, 0.186253860032785, 27.59128196499889, 248.32153768499, 0.01211767179672003, 0.01211767179672003),
    @('sCs', 'Cp', 'Slc40a1', 'ECs', 3, 1, 21.47235233333333, 64.417057, 0.3032756028926696, 0.3032756028926696, 2, 0.6666666666666666, 0.1166653333333333, 0.349996, 0.003627699309917676, 0.003627699309917676, 2.505079142419111, 22.545712281772, 0.001100192695328605, 0.001100192695328605),
    @('sCs', 'Cp', 'Slc40a1', 'FAPs', 3, 1, 21.47235233333333, 64.417057, 0.3032756028926696, 0.3032756028926696, 3, 1, 12.00209833333333, 36.00629499999999, 0.3732042981182421, 0.3732042981182421, 257.7132841526461, 2319.419557373815, 0.1131837585139455, 0.1131837585139455),
    @('sCs', 'Cp', 'Slc40a1', 'M2', 3, 1, 21.47235233333333, 64.417057, 0.3032756028926696, 0.3032756028926696, 3, 1, 14.050981, 42.15294300000001, 0.4369141425390552, 0.4369141425390553, 301.7076146609724, 2715.368531948751, 0.1325053999908658, 0.1325053999908657),
    @('sCs', 'Cp', 'Slc40a1', 'sCs', 3, 1, 21.47235233333333, 64.417057, 0.3032756028926696, 0.3032756028926696, 3, 1, 5.989848333333334, 17.969545, 0.186253860032785, 0.186253860032785, 128.6161338365628, 1157.545204529065, 0.05648625169252979, 0.05648625169252978)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
